$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Calculation options: enable iterative calculation (iterateDelta) ---
$excel.Iteration = $true
$excel.MaxChange = 0.0001

# --- Update shared text (Area responsable / Nota) ---
$ws.Range("I8").Value = "Departamento de Inventarios  (UPP) "
$ws.Range("L8").Value = "En el periodo que se informa no se realizo ningun inventario de bajas practicadas a bienes muebles."

# --- Update reporting period dates (row 8): 2022-01-01/2022-06-30 -> 2022-07-01/2022-12-31 ---
$ws.Range("B8").Value = 44743
$ws.Range("C8").Value = 44926

# --- Update validation/update dates: 2022-07-11 -> 2023-01-10 ---
$ws.Range("J8").Value = 44936
$ws.Range("K8").Value = 44936

# --- Column L width ---
$ws.Columns(12).ColumnWidth = 34.3

# --- Row heights ---
$ws.Rows(3).RowHeight = 27.75
$ws.Rows(8).RowHeight = 47.25

# --- Selection moved to B14 ---
$ws.Range("B14").Select()

# --- Page setup: portrait orientation ---
$ws.PageSetup.Orientation = 1

Write-Output "done"
